# Reorder the last two slides of the deck: the slide that is currently in
# position 30 ("raven.gif" placeholder slide) moves to position 29, and the
# slide that is currently in position 29 ("Generalized Polynomial Chaos")
# shifts down to position 30.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(30)
$s.MoveTo(29)
